$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.603177
$ws.Range("H2").Value = 19.809531
$ws.Range("I2").Value = 0.5135477412645301
$ws.Range("J2").Value = 0.5135477412645302
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 58.244295596578
$ws.Range("R2").Value = 524.198660369202
$ws.Range("S2").Value = 0.03294484077335207
$ws.Range("T2").Value = 0.03294484077335207
$ws.Range("G3").Value = 6.603177
$ws.Range("H3").Value = 19.809531
$ws.Range("I3").Value = 0.5135477412645301
$ws.Range("J3").Value = 0.5135477412645302
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 361.269949546695
$ws.Range("R3").Value = 3251.429545920255
$ws.Range("S3").Value = 0.2043458649830778
$ws.Range("T3").Value = 0.2043458649830778
$ws.Range("G4").Value = 6.603177
$ws.Range("H4").Value = 19.809531
$ws.Range("I4").Value = 0.5135477412645301
$ws.Range("J4").Value = 0.5135477412645302
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 144.66348023491
$ws.Range("R4").Value = 1301.97132211419
$ws.Range("S4").Value = 0.08182630201365304
$ws.Range("T4").Value = 0.08182630201365307
$ws.Range("G5").Value = 6.603177
$ws.Range("H5").Value = 19.809531
$ws.Range("I5").Value = 0.5135477412645301
$ws.Range("J5").Value = 0.5135477412645302
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 343.7406540410429
$ws.Range("R5").Value = 3093.665886369386
$ws.Range("S5").Value = 0.1944307334944471
$ws.Range("T5").Value = 0.1944307334944472
$ws.Range("I6").Value = 0.02944398858046029
$ws.Range("J6").Value = 0.0294439885804603
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 3.339405933710889
$ws.Range("R6").Value = 30.054653403398
$ws.Range("S6").Value = 0.001888875050111449
$ws.Range("T6").Value = 0.001888875050111449
$ws.Range("I7").Value = 0.02944398858046029
$ws.Range("J7").Value = 0.0294439885804603
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.01171606226951891
$ws.Range("T7").Value = 0.01171606226951892
$ws.Range("I8").Value = 0.02944398858046029
$ws.Range("J8").Value = 0.0294439885804603
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 8.294204253645555
$ws.Range("R8").Value = 74.64783828281
$ws.Range("S8").Value = 0.004691467819795669
$ws.Range("T8").Value = 0.00469146781979567
$ws.Range("I9").Value = 0.02944398858046029
$ws.Range("J9").Value = 0.0294439885804603
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 19.70818889652367
$ws.Range("R9").Value = 177.373700068713
$ws.Range("S9").Value = 0.01114758344103426
$ws.Range("T9").Value = 0.01114758344103426
$ws.Range("G10").Value = 3.441487333333333
$ws.Range("H10").Value = 10.324462
$ws.Range("I10").Value = 0.2676541983690312
$ws.Range("J10").Value = 0.2676541983690313
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 30.35614606946711
$ws.Range("R10").Value = 273.205314625204
$ws.Range("S10").Value = 0.01717040936812305
$ws.Range("T10").Value = 0.01717040936812306
$ws.Range("G11").Value = 3.441487333333333
$ws.Range("H11").Value = 10.324462
$ws.Range("I11").Value = 0.2676541983690312
$ws.Range("J11").Value = 0.2676541983690313
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 188.2890546897233
$ws.Range("R11").Value = 1694.60149220751
$ws.Range("S11").Value = 0.1065023254651974
$ws.Range("T11").Value = 0.1065023254651974
$ws.Range("G12").Value = 3.441487333333333
$ws.Range("H12").Value = 10.324462
$ws.Range("I12").Value = 0.2676541983690312
$ws.Range("J12").Value = 0.2676541983690313
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 75.39666660826445
$ws.Range("R12").Value = 678.5699994743801
$ws.Range("S12").Value = 0.04264677168482608
$ws.Range("T12").Value = 0.04264677168482608
$ws.Range("G13").Value = 3.441487333333333
$ws.Range("H13").Value = 10.324462
$ws.Range("I13").Value = 0.2676541983690312
$ws.Range("J13").Value = 0.2676541983690313
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 179.1530208616193
$ws.Range("R13").Value = 1612.377187754574
$ws.Range("S13").Value = 0.1013346918508846
$ws.Range("T13").Value = 0.1013346918508847
$ws.Range("G14").Value = 2.434707333333333
$ws.Range("H14").Value = 7.304122
$ws.Range("I14").Value = 0.1893540717859783
$ws.Range("J14").Value = 0.1893540717859783
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 21.47569474721378
$ws.Range("R14").Value = 193.281252724924
$ws.Range("S14").Value = 0.01214734141253207
$ws.Range("T14").Value = 0.01214734141253207
$ws.Range("G15").Value = 2.434707333333333
$ws.Range("H15").Value = 7.304122
$ws.Range("I15").Value = 0.1893540717859783
$ws.Range("J15").Value = 0.1893540717859783
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 133.2065754824233
$ws.Range("R15").Value = 1198.85917934181
$ws.Range("S15").Value = 0.07534590940249561
$ws.Range("T15").Value = 0.07534590940249562
$ws.Range("G16").Value = 2.434707333333333
$ws.Range("H16").Value = 7.304122
$ws.Range("I16").Value = 0.1893540717859783
$ws.Range("J16").Value = 0.1893540717859783
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 53.33996592753111
$ws.Range("R16").Value = 480.0596933477801
$ws.Range("S16").Value = 0.03017079469052385
$ws.Range("T16").Value = 0.03017079469052385
$ws.Range("G17").Value = 2.434707333333333
$ws.Range("H17").Value = 7.304122
$ws.Range("I17").Value = 0.1893540717859783
$ws.Range("J17").Value = 0.1893540717859783
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 126.7432163575993
$ws.Range("R17").Value = 1140.688947218394
$ws.Range("S17").Value = 0.07169002628042674
$ws.Range("T17").Value = 0.07169002628042677
